# Updated cryptos list (price/volume refresh) on Fri Jan 26 04:30:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.163.36"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.225.57"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'294.10"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("D6").Value = "'88.09"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").Value = "'30.75"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").Value = "'50.98"
$ws.Range("E11").Value = "  +6.47%  "

$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("D14").Value = "'6.44"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'13.85"
$ws.Range("E15").Value = "  -1.18%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.205.75"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.738"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").Value = "2.037.71"
$ws.Range("E18").Value = "  -20.37%  "

$ws.Range("D19").Value = "40.075.02"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").Value = "'11.27"
$ws.Range("E21").Value = "  -4.70%  "

$ws.Range("D22").Value = "'5.79"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").Value = "'65.74"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'236.27"
$ws.Range("E24").Value = "  +0.59%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("D26").Value = "'2.48"
$ws.Range("E26").Value = "  +0.79%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").Value = "'23.25"
$ws.Range("E28").Value = "  +2.76%  "

$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("E30").Value = "  -10.34%  "

$ws.Range("D31").Value = "'158.94"
$ws.Range("E31").Value = "  +3.41%  "

$ws.Range("D32").Value = "'31.94"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("E35").Value = "  +6.18%  "

$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("D38").Value = "'0.113"
$ws.Range("E38").Value = "  +1.64%  "

$ws.Range("E39").Value = "  +3.14%  "

$ws.Range("D40").Value = "'0.0994"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").Value = "'15.62"
$ws.Range("E41").Value = "  -1.81%  "

$ws.Range("D42").Value = "2.088.54"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("E43").Value = "  -2.99%  "

$ws.Range("D44").Value = "'19.29"
$ws.Range("E44").Value = "  +8.71%  "

$ws.Range("D45").Value = "'10.13"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("D47").Value = "'2.75"
$ws.Range("E47").Value = "  +2.69%  "

$ws.Range("E48").Value = "  -13.31%  "

$ws.Range("D49").Value = "2.441.47"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("E50").Value = "  +1.98%  "

$ws.Range("E51").Value = "  +3.76%  "
